$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 301. Everything at/after row 301
# (Region de O'Higgins durazno records) shifts down by one, and the
# sheet dimension grows from A1:T378 to A1:T379.
$ws.Rows.Item(301).Insert()

# Populate the newly inserted row 301 with the new price record
# (Vega Monumental Concepcion - Durazno - Early Majestic - Primera).
$ws.Range("A301").Value2 = 11
$ws.Range("B301").Value2 = "Vega Monumental Concepción"
$ws.Range("C301").Value2 = "Bíobío"
$ws.Range("D301").Value2 = 45275
$ws.Range("E301").Value2 = 8
$ws.Range("F301").Value2 = "Fruta"
$ws.Range("G301").Value2 = 100103
$ws.Range("H301").Value2 = "Frutos de hueso (carozo)"
$ws.Range("I301").Value2 = 100103004
$ws.Range("J301").Value2 = "Durazno"
$ws.Range("K301").Value2 = "Early Majestic"
$ws.Range("L301").Value2 = "Primera"
$ws.Range("M301").Value2 = 100
$ws.Range("N301").Value2 = 15000
$ws.Range("O301").Value2 = 16000
$ws.Range("P301").Value2 = 15500
$ws.Range("Q301").Value2 = "$/caja 15 kilos empedrada"
$ws.Range("R301").Value2 = "Región de O'Higgins"
$ws.Range("S301").Value2 = 1033
$ws.Range("T301").Value2 = 15
